# Regenerated save_data: column G ("K") values were recalculated from the
# (re-generated) source simulation and are being rewritten here with their
# new values. Column G is the 7th column on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 2
    8  = 2
    9  = 2
    10 = 2
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    20 = 2
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 2
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 1
    33 = 2
    34 = 0
    35 = 2
    36 = 3
    37 = 1
    38 = 3
    39 = 0
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 0
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 0
    51 = 1
    52 = 0
    53 = 1
    54 = 0
    55 = 3
    56 = 1
    57 = 2
    58 = 0
    59 = 1
    60 = 1
    61 = 3
    62 = 2
    63 = 2
    64 = 0
    65 = 1
    66 = 1
    68 = 2
    69 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
